$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 4 is the "Name" property; set its value cell (B4) which was previously empty.
$ws.Range("B4").Value = "GenreactiviteVs"

# Row 8 is the "Date" property; update the generation timestamp.
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
